$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 19): "2356. Number of Unique Subjects Taught by Each Teacher"
$ws.Cells.Item(19, 1).Value = "2356. Number of Unique Subjects Taught by Each Teacher"
$ws.Cells.Item(19, 2).Value = "Easy"
$ws.Cells.Item(19, 3).Value = "Sorting and Grouping"
$ws.Cells.Item(19, 4).Value = "Use COUNT(DISTINCT subject_id) and GROUP BY teacher_id"
$ws.Cells.Item(19, 5).Value = "https://leetcode.com/problems/number-of-unique-subjects-taught-by-each-teacher/solutions/3871227/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "

# Match formatting used by the other rows: green "Easy" fill on column B
$ws.Range("B19").Interior.Color = $ws.Range("B2").Interior.Color()

# Turn the link text into a real, clickable hyperlink, same as the rest of the table
$ws.Hyperlinks.Add($ws.Range("E19"), "https://leetcode.com/problems/number-of-unique-subjects-taught-by-each-teacher/solutions/3871227/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 ")

# Re-apply the same "Hyperlink" cell style used by the other link cells in column E
$ws.Range("E19").Style = "Hyperlink"

# Grow the table so it covers the new row too
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E19"))

# Update the active cell selection like in the edited workbook
$ws.Range("E27").Select()
